$d = $word.ActiveDocument

# 0. Cosmetic rename of the two inline pictures' internal docPr/cNvPr
#    "name" id (envelope.png: image5.png -> image4.png; linkedin.png:
#    image2.png -> image3.png). Word's real InlineShape object does not
#    expose a writable Name property (only floating Shape.Name is
#    settable), so this is attempted best-effort and simply has no
#    effect if unsupported by the host - it does not affect the rest of
#    the script.
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shp = $d.InlineShapes.Item($i)
    if ($shp.AlternativeText -eq "envelope.png") {
        try { $shp.Name = "image4.png" } catch { }
    } elseif ($shp.AlternativeText -eq "linkedin.png") {
        try { $shp.Name = "image3.png" } catch { }
    }
}

# 1. Qualifications section: HTML -> HTML5
$d.Content.Find.Execute(
    "Experienced in JavaScript, Angular, jQuery, CSS, HTML",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Experienced in JavaScript, Angular, jQuery, CSS, HTML5",
    2)

# 2. Qualifications section: drop trailing ", Node.js"
$d.Content.Find.Execute(
    "Proficient in C#, PHP, SQL, MVC, Python, Node.js",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Proficient in C#, PHP, SQL, MVC, Python",
    2)

# 3. Qualifications section: append ", Node.js" to the exposure line
$d.Content.Find.Execute(
    "Exposure to Vue.js, React, AWS, IIS, Flask",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Exposure to Vue.js, React, AWS, IIS, Flask, Node.js",
    2)

# 4. Experience bullet: add React to the stack list
$d.Content.Find.Execute(
    "ull stack apps leveraging: AngularJS, C3.js, D3.js, Sass, C# .NET, SQL Server",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "ull stack apps leveraging: AngularJS, React, C3.js, D3.js, Sass, C# .NET, SQL Server",
    2)

# 5. Education bullet: mention the Oracle PL/SQL project alongside the PHP clone
$d.Content.Find.Execute(
    "Built an Amazon.com like clone for a final software engineering project using PHP",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Oracle PL/SQL client/server final project and an Amazon.com like clone for a final software engineering project using PHP",
    2)
